# Natmi following Dr Hou advice
# Update the Dkk2-Lrp5 LR-pair results for the 3 target clusters (rows 2-4)
# after recomputing the number of ligand/receptor expressing cells (1 -> 3)
# and the downstream average/total expression & specificity metrics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster: ECs)
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 8.398553333333334
$ws.Range("H2").Value = 25.19566
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 6.215523666666666
$ws.Range("N2").Value = 18.646571
$ws.Range("O2").Value = 0.2852115546146347
$ws.Range("P2").Value = 0.2852115546146347
$ws.Range("Q2").Value = 52.20140700909555
$ws.Range("R2").Value = 469.81266308186
$ws.Range("S2").Value = 0.2852115546146347
$ws.Range("T2").Value = 0.2852115546146347

# Row 3 (Target cluster: FAPs)
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 8.398553333333334
$ws.Range("H3").Value = 25.19566
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 9.484070666666666
$ws.Range("N3").Value = 28.452212
$ws.Range("O3").Value = 0.4351952762116512
$ws.Range("P3").Value = 0.4351952762116512
$ws.Range("Q3").Value = 79.65247331110223
$ws.Range("R3").Value = 716.87225979992
$ws.Range("S3").Value = 0.4351952762116512
$ws.Range("T3").Value = 0.4351952762116512

# Row 4 (Target cluster: sCs)
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 8.398553333333334
$ws.Range("H4").Value = 25.19566
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.093084
$ws.Range("N4").Value = 18.279252
$ws.Range("O4").Value = 0.2795931691737141
$ws.Range("P4").Value = 0.2795931691737141
$ws.Range("Q4").Value = 51.17309093848001
$ws.Range("R4").Value = 460.55781844632
$ws.Range("S4").Value = 0.2795931691737141
$ws.Range("T4").Value = 0.2795931691737141
